$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateValue = 43221

# Dates: fill A99:A101 and reuse the existing short-date style from A98
# (keeps cellXfs untouched instead of generating a brand-new custom format)
$ws.Range("A99:A101").Value = $dateValue
$ws.Range("A98").Copy()
$ws.Range("A99:A101").PasteSpecial(-4122)

# Row 99: description text is entered before the image filename
$ws.Cells.Item(99, 3).Value = "Debuggen von diesen Großen konstrukten mit allen einzelteilen. Nur keine Panik!"
$ws.Cells.Item(99, 2).Value = "2018-05-01 8.jpg"

# Row 100: image filename entered before the description text
$ws.Cells.Item(100, 2).Value = "2018-05-01 9.jpg"
$ws.Cells.Item(100, 3).Value = "Wenn mich die Hardware verlässt, kann die Zeit wenigstens für die Software verwendet werden. Es gibt jetzt eine stabile, relativ schnelle und zuverlässige zwei-Wege Kommunikation zwischen Basis und Arduino (Drohne)"

# Row 101: image filename entered before the description text
$ws.Cells.Item(101, 2).Value = "2018-05-01 10.jpg"
$ws.Cells.Item(101, 3).Value = "Die neuen ESCs sind Bestellt und garantiert bis Freitag angekommen und einsatzbereit"

$ws.Range("A79").Select()
$ws.Cells.Item(101, 3).Select()
